# Applies scheduled-runner price/profit updates to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5498.143
$ws.Range("I33").Value = 7314.3335
$ws.Range("K33").Value = 7314.3335
$ws.Range("M33").Value = -7085.3335
$ws.Range("H132").Value = 5679.0625
$ws.Range("I132").Value = 4153.654
$ws.Range("K132").Value = 12460.962
$ws.Range("M132").Value = -9930.962000000001
$ws.Range("H138").Value = 3795.318
$ws.Range("I138").Value = 5938.8423
$ws.Range("J138").Value = 3205.0725
$ws.Range("K138").Value = 17816.5269
$ws.Range("L138").Value = 9615.217500000001
$ws.Range("M138").Value = -12676.5269
$ws.Range("N138").Value = -19895.2175

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 4817.778
$ws.Range("I31").Value = 4817.778
$ws.Range("K31").Value = 4817.778
$ws.Range("M31").Value = -4523.778
$ws.Range("H32").Value = 28232.867
$ws.Range("I32").Value = 26500
$ws.Range("J32").Value = 28499.46
$ws.Range("K32").Value = 26500
$ws.Range("L32").Value = 28499.46
$ws.Range("M32").Value = -26213
$ws.Range("N32").Value = -29073.46
$ws.Range("H74").Value = 43902.25
$ws.Range("I74").Value = 52053.2
$ws.Range("J74").Value = 3147.5
$ws.Range("K74").Value = 52053.2
$ws.Range("L74").Value = 3147.5
$ws.Range("M74").Value = -51179.2
$ws.Range("N74").Value = -4895.5
$ws.Range("H77").Value = 43902.25
$ws.Range("I77").Value = 52053.2
$ws.Range("J77").Value = 3147.5
$ws.Range("K77").Value = 260266
$ws.Range("L77").Value = 15737.5
$ws.Range("M77").Value = -255898
$ws.Range("N77").Value = -24473.5
$ws.Range("H88").Value = 2505.6365
$ws.Range("I88").Value = 1997
$ws.Range("J88").Value = 2618.6667
$ws.Range("K88").Value = 1997
$ws.Range("L88").Value = 2618.6667
$ws.Range("M88").Value = -1591
$ws.Range("N88").Value = -3430.6667
$ws.Range("H91").Value = 2505.6365
$ws.Range("I91").Value = 1997
$ws.Range("J91").Value = 2618.6667
$ws.Range("K91").Value = 1997
$ws.Range("L91").Value = 2618.6667
$ws.Range("M91").Value = -593
$ws.Range("N91").Value = -5426.6667
$ws.Range("H132").Value = 31920.893
$ws.Range("I132").Value = 35696.344
$ws.Range("J132").Value = 7758
$ws.Range("K132").Value = 107089.032
$ws.Range("L132").Value = 23274
$ws.Range("M132").Value = -104559.032
$ws.Range("N132").Value = -28334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H94").Value = 1647.4
$ws.Range("I94").Value = 1059.1428
$ws.Range("J94").Value = 2162.125
$ws.Range("K94").Value = 1059.1428
$ws.Range("L94").Value = 2162.125
$ws.Range("M94").Value = -608.1428000000001
$ws.Range("N94").Value = -3064.125
$ws.Range("H107").Value = 4199.5
$ws.Range("I107").Value = 4199.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4199.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2279.5
$ws.Range("N107").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 20996.834
$ws.Range("J80").Value = 20996.834
$ws.Range("L80").Value = 20996.834
$ws.Range("N80").Value = -23242.834
$ws.Range("H83").Value = 20996.834
$ws.Range("J83").Value = 20996.834
$ws.Range("L83").Value = 62990.50199999999
$ws.Range("N83").Value = -74222.50199999999
$ws.Range("H105").Value = 2049.625
$ws.Range("I105").Value = 1579.4
$ws.Range("J105").Value = 2833.3333
$ws.Range("K105").Value = 1579.4
$ws.Range("L105").Value = 2833.3333
$ws.Range("M105").Value = 167.5999999999999
$ws.Range("N105").Value = -6327.3333
$ws.Range("H129").Value = 79999
$ws.Range("J129").Value = 79999
$ws.Range("L129").Value = 79999
$ws.Range("N129").Value = -89999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 6699.75
$ws.Range("I58").Value = 5266.6665
$ws.Range("K58").Value = 15799.9995
$ws.Range("M58").Value = -15671.9995
$ws.Range("H64").Value = 5088.875
$ws.Range("I64").Value = 2237.3333
$ws.Range("J64").Value = 6799.8
$ws.Range("K64").Value = 6711.999899999999
$ws.Range("L64").Value = 20399.4
$ws.Range("M64").Value = -6441.999899999999
$ws.Range("N64").Value = -20939.4
$ws.Range("H67").Value = 5088.875
$ws.Range("I67").Value = 2237.3333
$ws.Range("J67").Value = 6799.8
$ws.Range("K67").Value = 6711.999899999999
$ws.Range("L67").Value = 20399.4
$ws.Range("M67").Value = -5775.999899999999
$ws.Range("N67").Value = -22271.4
$ws.Range("H82").Value = 6006.5
$ws.Range("I82").Value = 2013
$ws.Range("K82").Value = 6039
$ws.Range("M82").Value = -5633
$ws.Range("H85").Value = 6006.5
$ws.Range("I85").Value = 2013
$ws.Range("K85").Value = 6039
$ws.Range("M85").Value = -4635
$ws.Range("H131").Value = 7711001
$ws.Range("I131").Value = 1521
$ws.Range("J131").Value = 9112725
$ws.Range("K131").Value = 4563
$ws.Range("L131").Value = 27338175
$ws.Range("M131").Value = 477
$ws.Range("N131").Value = -27348255

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 37371.75
$ws.Range("J15").Value = 46662.332
$ws.Range("L15").Value = 46662.332
$ws.Range("N15").Value = -47238.332
$ws.Range("H81").Value = 37371.75
$ws.Range("J81").Value = 46662.332
$ws.Range("L81").Value = 46662.332
$ws.Range("N81").Value = -48658.332
$ws.Range("H84").Value = 37371.75
$ws.Range("J84").Value = 46662.332
$ws.Range("L84").Value = 139986.996
$ws.Range("N84").Value = -149970.996
$ws.Range("H132").Value = 77749.84
$ws.Range("I132").Value = 84152.836
$ws.Range("J132").Value = 914
$ws.Range("K132").Value = 252458.508
$ws.Range("L132").Value = 2742
$ws.Range("M132").Value = -249928.508
$ws.Range("N132").Value = -7802
$ws.Range("H139").Value = 59997
$ws.Range("J139").Value = 59997
$ws.Range("L139").Value = 59997
$ws.Range("N139").Value = -70277

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 445174
$ws.Range("I7").Value = 537636.9399999999
$ws.Range("J7").Value = 5975
$ws.Range("K7").Value = 537636.9399999999
$ws.Range("L7").Value = 5975
$ws.Range("M7").Value = -537524.9399999999
$ws.Range("N7").Value = -6199
$ws.Range("H42").Value = 32005.25
$ws.Range("J42").Value = 39998.5
$ws.Range("L42").Value = 39998.5
$ws.Range("N42").Value = -41124.5
$ws.Range("H49").Value = 32005.25
$ws.Range("J49").Value = 39998.5
$ws.Range("L49").Value = 39998.5
$ws.Range("N49").Value = -40292.5
$ws.Range("H82").Value = 3081.1765
$ws.Range("I82").Value = 1968.875
$ws.Range("K82").Value = 1968.875
$ws.Range("M82").Value = -1607.875
$ws.Range("H85").Value = 3081.1765
$ws.Range("I85").Value = 1968.875
$ws.Range("K85").Value = 1968.875
$ws.Range("M85").Value = -720.875
$ws.Range("H126").Value = 445174
$ws.Range("I126").Value = 537636.9399999999
$ws.Range("J126").Value = 5975
$ws.Range("K126").Value = 1612910.82
$ws.Range("L126").Value = 17925
$ws.Range("M126").Value = -1610440.82
$ws.Range("N126").Value = -22865
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 914197.5600000001
$ws.Range("I81").Value = 2083.625
$ws.Range("J81").Value = 3346501.2
$ws.Range("K81").Value = 4167.25
$ws.Range("L81").Value = 6693002.4
$ws.Range("M81").Value = -3106.25
$ws.Range("N81").Value = -6695124.4
$ws.Range("H84").Value = 914197.5600000001
$ws.Range("I84").Value = 2083.625
$ws.Range("J84").Value = 3346501.2
$ws.Range("K84").Value = 20836.25
$ws.Range("L84").Value = 33465012
$ws.Range("M84").Value = -15532.25
$ws.Range("N84").Value = -33475620
$ws.Range("H93").Value = 252694.5
$ws.Range("J93").Value = 252694.5
$ws.Range("L93").Value = 252694.5
$ws.Range("N93").Value = -257686.5
$ws.Range("H113").Value = 1328.9333
$ws.Range("I113").Value = 1157.8636
$ws.Range("K113").Value = 3473.5908
$ws.Range("M113").Value = -1303.5908
$ws.Range("H136").Value = 2568.1064
$ws.Range("I136").Value = 2265.932
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 6797.795999999999
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -4247.795999999999
$ws.Range("N136").Value = -26100
